$d = $word.ActiveDocument
$full = $d.Content

# ---------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------

# Locate the n-th (1-indexed) occurrence of $text in the document and
# return a Range covering it (collapsed search, left-to-right).
function Find-Nth($text, $n) {
    $rng = $d.Content
    $rng.Start = 0
    for ($i = 0; $i -lt $n; $i++) {
        $found = $rng.Find.Execute($text, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if (-not $found) { return $null }
        if ($i -lt $n - 1) {
            $endPos = $rng.End
            $rng.SetRange($endPos, $full.End)
        }
    }
    return $rng
}

# Toggling a character-level property on/off forces the COM host to
# materialize the range as its own <w:r> instead of silently merging it
# back into an adjacent run with identical formatting.
function Force-Split($rng) {
    $rng.Bold = 1
    $rng.Bold = 0
}

# Insert literal $text right before document position $pos, as its own
# run (split off from whatever run used to occupy that boundary).
# Returns the Range covering the newly-inserted text.
function Insert-SplitRun($pos, $text) {
    $insPoint = $d.Range($pos, $pos)
    $insPoint.InsertBefore($text)
    $newRng = $d.Range($pos, $pos + $text.Length)
    Force-Split $newRng
    return $newRng
}

# ---------------------------------------------------------------------
# The template has two copies of the Name/Strand/Department/Age/Gender
# form row. The first copy already reads ${pn} / ${pc} / ${dept} /
# ${age} / ${g}. The second copy (further down, preceding the readings
# table) still reads the bare {pn} / {pc} / {dept} / {age} / {g} and is
# the one this change updates to the ${...} placeholder syntax as well.
# ---------------------------------------------------------------------

# --- {pn} -> ${pn} (2nd occurrence in the document) ---
$pn = Find-Nth "{pn}" 2
Insert-SplitRun $pn.Start "$" | Out-Null

# --- " {pc}" -> "${pc}" (only 1 occurrence; the first table has no
#     leading space before its "{pc}") ---
$pc = Find-Nth " {pc}" 1
# Drop the leading space.
$spaceRng = $d.Range($pc.Start, $pc.Start + 1)
$spaceRng.Delete()
# The remaining "{" character (now at $pc.Start) becomes "$".
$braceRng = $d.Range($pc.Start, $pc.Start + 1)
$braceRng.Text = "$"
Force-Split ($d.Range($pc.Start, $pc.Start + 1))
# Re-insert "{" as its own run right after the new "$" run.
Insert-SplitRun ($pc.Start + 1) "{" | Out-Null
# Keep "pc" / "}" as separate runs (they were before the edit too).
Force-Split ($d.Range($pc.Start + 2, $pc.Start + 4))

# --- {dept} -> ${dept} (2nd occurrence) ---
$dept = Find-Nth "{dept}" 2
Insert-SplitRun $dept.Start "$" | Out-Null

# --- " {age}" -> " ${age}" (only 1 occurrence) ---
$age = Find-Nth " {age}" 1
Insert-SplitRun ($age.Start + 1) "$" | Out-Null

# --- {g} -> ${g} (2nd occurrence) ---
$g = Find-Nth "{g}" 2
Insert-SplitRun $g.Start "$" | Out-Null
